$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")
$newTumorQuery = @'
MATCH (s:study)<--(p:participant)<--(samp:sample)
WHERE s.study_name in ["Human Tumor Atlas Network (HTAN) primary sequencing data"]
WITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor
RETURN  
 coalesce(samp.sample_id, '') as `Sample ID`,
 coalesce(p.participant_id,'') as `Participant ID`,
 coalesce(s.study_name, '') as `Study Name`,
 coalesce(s.phs_accession,'') as `Accession`,
 coalesce(samp.sample_tumor_status,'') as `Tumor`,
coalesce(samp.sample_type,'') as `Analyte Type`
  ORDER By samp.sample_id LIMIT 100
'@
$ws.Range("B3").Value = $newTumorQuery
